$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) cells are treated as text so values like
# "9.80", "1.00", "0.999" keep their exact formatting instead of being
# reinterpreted as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.816.61"
$ws.Range("E2").Value = "  +6.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.626.97"
$ws.Range("E3").Value = "  +6.07%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.39"
$ws.Range("E5").Value = "  +4.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.85"
$ws.Range("E6").Value = "  +7.49%  "

$ws.Range("E7").Value = "  +3.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.619.70"
$ws.Range("E8").Value = "  +6.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +2.89%  "

$ws.Range("E11").Value = "  +4.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.34"
$ws.Range("E12").Value = "  +6.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  +6.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.80"
$ws.Range("E14").Value = "  +5.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.205.61"
$ws.Range("E15").Value = "  +6.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.629.03"
$ws.Range("E16").Value = "  +5.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.40"
$ws.Range("E17").Value = "  +5.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.657.49"
$ws.Range("E18").Value = "  +6.97%  "

$ws.Range("E19").Value = "  +5.07%  "

$ws.Range("E20").Value = "  +0.53%  "

$ws.Range("E21").Value = "  +4.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.49"
$ws.Range("E22").Value = "  +6.25%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.33"
$ws.Range("E23").Value = "  +17.69%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.41"
$ws.Range("E24").Value = "  +8.90%  "

$ws.Range("E25").Value = "  +8.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.14"
$ws.Range("E26").Value = "  +1.33%  "

$ws.Range("E27").Value = "  +6.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.26"
$ws.Range("E28").Value = "  +4.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("E29").Value = "  +7.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.41"
$ws.Range("E30").Value = "  +3.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("E31").Value = "  +11.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.25"
$ws.Range("E32").Value = "  +5.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "621.37"
$ws.Range("E33").Value = "  +6.41%  "

$ws.Range("E34").Value = "  +8.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.24"
$ws.Range("E35").Value = "  +4.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0830"
$ws.Range("E36").Value = "  +8.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.413"
$ws.Range("E37").Value = "  +8.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.27"
$ws.Range("E38").Value = "  +4.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.327.57"
$ws.Range("E42").Value = "  +6.37%  "

$ws.Range("E43").Value = "  +4.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0448"

$ws.Range("E45").Value = "  +7.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").Value = "  +5.69%  "

$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("E48").Value = "  +7.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("E49").Value = "  +3.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.31"
$ws.Range("E50").Value = "  +4.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.05%  "
